$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '69.345.28'
$ws.Range("E2").Value = '  +4.71%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.498.81'

$ws.Range("E4").Value = '  +0.06%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '587.25'
$ws.Range("E5").Value = '  +2.22%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '182.46'
$ws.Range("E6").Value = '  +7.32%  '

$ws.Range("B7").Value = 'USDC'
$ws.Range("C7").Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.00'
$ws.Range("E7").Value = '  +0.13%  '

$ws.Range("B8").Value = 'LidoStakedEther'
$ws.Range("C8").Value = 'https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '3.495.01'
$ws.Range("E8").Value = '  +13.71%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.531'
$ws.Range("E9").Value = '  +4.10%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '6.56'
$ws.Range("E10").Value = '  +3.74%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.157'
$ws.Range("E11").Value = '  +5.43%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.487'
$ws.Range("E12").Value = '  +3.71%  '

$ws.Range("B13").Value = 'ShibaInu'
$ws.Range("C13").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000249'
$ws.Range("E13").Value = '  +4.17%  '

$ws.Range("B14").Value = 'Avalanche'
$ws.Range("C14").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '38.24'
$ws.Range("E14").Value = '  +6.47%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '4.081.80'
$ws.Range("E15").Value = '  +13.73%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '69.492.39'
$ws.Range("E16").Value = '  +5.04%  '

$ws.Range("E17").Value = '  +1.24%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.491.89'
$ws.Range("E18").Value = '  +13.55%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.39'
$ws.Range("E19").Value = '  +6.31%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '16.79'
$ws.Range("E20").Value = '  +1.74%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '503.47'
$ws.Range("E21").Value = '  +4.10%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '8.91'
$ws.Range("E22").Value = '  +16.32%  '

$ws.Range("E23").Value = '  +6.52%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '86.31'
$ws.Range("E24").Value = '  +4.68%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '13.32'
$ws.Range("E25").Value = '  +5.44%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.37'
$ws.Range("E26").Value = '  +7.41%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.67'
$ws.Range("E27").Value = '  +4.57%  '

$ws.Range("E28").Value = '  -0.01%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.50'
$ws.Range("E29").Value = '  +11.68%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '8.06'
$ws.Range("E30").Value = '  +2.41%  '

$ws.Range("B31").Value = 'PEPE'
$ws.Range("C31").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.0000108'
$ws.Range("E31").Value = '  +19.99%  '

$ws.Range("B32").Value = 'PancakeSwap'
$ws.Range("C32").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.71'
$ws.Range("E32").Value = '  +4.35%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '30.54'
$ws.Range("E33").Value = '  +10.10%  '

$ws.Range("E34").Value = '  +5.15%  '

$ws.Range("E35").Value = '  +0.18%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '6.06'
$ws.Range("E36").Value = '  +9.00%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.01'
$ws.Range("E37").Value = '  +7.06%  '

$ws.Range("E38").Value = '  +10.18%  '

$ws.Range("B39").Value = 'Stacks'
$ws.Range("C39").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.10'
$ws.Range("E39").Value = '  +7.40%  '

$ws.Range("B40").Value = 'Arweave'
$ws.Range("C40").Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '46.54'
$ws.Range("E40").Value = '  -1.62%  '

$ws.Range("B41").Value = 'OKB'
$ws.Range("C41").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '50.24'
$ws.Range("E41").Value = '  +2.43%  '

$ws.Range("B42").Value = 'Kaspa'
$ws.Range("C42").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.127'
$ws.Range("E42").Value = '  +4.26%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '8.67'
$ws.Range("E43").Value = '  +5.41%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.78'
$ws.Range("E44").Value = '  +11.00%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.989.07'
$ws.Range("E45").Value = '  +7.65%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '404.85'
$ws.Range("E46").Value = '  +11.12%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0362'
$ws.Range("E47").Value = '  +5.42%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '27.56'
$ws.Range("E48").Value = '  +13.06%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '134.76'
$ws.Range("E49").Value = '  +0.09%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.43'
$ws.Range("E51").Value = '  +13.15%  '
